$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-7: 45243 -> 45244 (+1 day)
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 3).Value2 = 45244
}
